# Auto-generated edit script: append new sensor-log rows to match the
# authoritative commit diff (adds rows to ALERTS, PIR, Humidity,
# Temperature and Proximity sheets, each a contiguous block of new
# inline-string records for 2026-02-06).

$wb = $excel.ActiveWorkbook

# --- ALERTS: append 1 row(s) starting at row 10 ---
$ws = $wb.Worksheets.Item("ALERTS")

$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "2026-02-06"
$ws.Cells.Item(10, 2).Value = "09:55:50"
$ws.Cells.Item(10, 3).Value = "09:00"
$ws.Cells.Item(10, 4).Value = "Bathroom"
$ws.Cells.Item(10, 5).Value = "MINIMAL"
$ws.Cells.Item(10, 6).Value = "MINIMAL ALERT: Bathroom occupied, no motion > 20s."


# --- PIR: append 14 row(s) starting at row 178 ---
$ws = $wb.Worksheets.Item("PIR")

$ws.Cells.Item(178, 1).NumberFormat = "@"
$ws.Cells.Item(178, 1).Value = "2026-02-06"
$ws.Cells.Item(178, 2).Value = "09:54:51"
$ws.Cells.Item(178, 3).Value = "09:00"
$ws.Cells.Item(178, 4).Value = "Bathroom"
$ws.Cells.Item(178, 5).Value = "No Motion"
$ws.Cells.Item(178, 6).Value = "Inactive"

$ws.Cells.Item(179, 1).NumberFormat = "@"
$ws.Cells.Item(179, 1).Value = "2026-02-06"
$ws.Cells.Item(179, 2).Value = "09:54:51"
$ws.Cells.Item(179, 3).Value = "09:00"
$ws.Cells.Item(179, 4).Value = "Bathroom"
$ws.Cells.Item(179, 5).Value = "No Motion"
$ws.Cells.Item(179, 6).Value = "Inactive"

$ws.Cells.Item(180, 1).NumberFormat = "@"
$ws.Cells.Item(180, 1).Value = "2026-02-06"
$ws.Cells.Item(180, 2).Value = "09:54:56"
$ws.Cells.Item(180, 3).Value = "09:00"
$ws.Cells.Item(180, 4).Value = "Bathroom"
$ws.Cells.Item(180, 5).Value = "No Motion"
$ws.Cells.Item(180, 6).Value = "Inactive"

$ws.Cells.Item(181, 1).NumberFormat = "@"
$ws.Cells.Item(181, 1).Value = "2026-02-06"
$ws.Cells.Item(181, 2).Value = "09:55:01"
$ws.Cells.Item(181, 3).Value = "09:00"
$ws.Cells.Item(181, 4).Value = "Bathroom"
$ws.Cells.Item(181, 5).Value = "No Motion"
$ws.Cells.Item(181, 6).Value = "Inactive"

$ws.Cells.Item(182, 1).NumberFormat = "@"
$ws.Cells.Item(182, 1).Value = "2026-02-06"
$ws.Cells.Item(182, 2).Value = "09:55:06"
$ws.Cells.Item(182, 3).Value = "09:00"
$ws.Cells.Item(182, 4).Value = "Bathroom"
$ws.Cells.Item(182, 5).Value = "No Motion"
$ws.Cells.Item(182, 6).Value = "Inactive"

$ws.Cells.Item(183, 1).NumberFormat = "@"
$ws.Cells.Item(183, 1).Value = "2026-02-06"
$ws.Cells.Item(183, 2).Value = "09:55:08"
$ws.Cells.Item(183, 3).Value = "09:00"
$ws.Cells.Item(183, 4).Value = "Bathroom"
$ws.Cells.Item(183, 5).Value = "Motion Detected"
$ws.Cells.Item(183, 6).Value = "Active"

$ws.Cells.Item(184, 1).NumberFormat = "@"
$ws.Cells.Item(184, 1).Value = "2026-02-06"
$ws.Cells.Item(184, 2).Value = "09:55:15"
$ws.Cells.Item(184, 3).Value = "09:00"
$ws.Cells.Item(184, 4).Value = "Bathroom"
$ws.Cells.Item(184, 5).Value = "No Motion"
$ws.Cells.Item(184, 6).Value = "Inactive"

$ws.Cells.Item(185, 1).NumberFormat = "@"
$ws.Cells.Item(185, 1).Value = "2026-02-06"
$ws.Cells.Item(185, 2).Value = "09:55:16"
$ws.Cells.Item(185, 3).Value = "09:00"
$ws.Cells.Item(185, 4).Value = "Bathroom"
$ws.Cells.Item(185, 5).Value = "Motion Detected"
$ws.Cells.Item(185, 6).Value = "Active"

$ws.Cells.Item(186, 1).NumberFormat = "@"
$ws.Cells.Item(186, 1).Value = "2026-02-06"
$ws.Cells.Item(186, 2).Value = "09:55:24"
$ws.Cells.Item(186, 3).Value = "09:00"
$ws.Cells.Item(186, 4).Value = "Bathroom"
$ws.Cells.Item(186, 5).Value = "No Motion"
$ws.Cells.Item(186, 6).Value = "Inactive"

$ws.Cells.Item(187, 1).NumberFormat = "@"
$ws.Cells.Item(187, 1).Value = "2026-02-06"
$ws.Cells.Item(187, 2).Value = "09:55:27"
$ws.Cells.Item(187, 3).Value = "09:00"
$ws.Cells.Item(187, 4).Value = "Bathroom"
$ws.Cells.Item(187, 5).Value = "Motion Detected"
$ws.Cells.Item(187, 6).Value = "Active"

$ws.Cells.Item(188, 1).NumberFormat = "@"
$ws.Cells.Item(188, 1).Value = "2026-02-06"
$ws.Cells.Item(188, 2).Value = "09:55:34"
$ws.Cells.Item(188, 3).Value = "09:00"
$ws.Cells.Item(188, 4).Value = "Bathroom"
$ws.Cells.Item(188, 5).Value = "No Motion"
$ws.Cells.Item(188, 6).Value = "Inactive"

$ws.Cells.Item(189, 1).NumberFormat = "@"
$ws.Cells.Item(189, 1).Value = "2026-02-06"
$ws.Cells.Item(189, 2).Value = "09:55:39"
$ws.Cells.Item(189, 3).Value = "09:00"
$ws.Cells.Item(189, 4).Value = "Bathroom"
$ws.Cells.Item(189, 5).Value = "No Motion"
$ws.Cells.Item(189, 6).Value = "Inactive"

$ws.Cells.Item(190, 1).NumberFormat = "@"
$ws.Cells.Item(190, 1).Value = "2026-02-06"
$ws.Cells.Item(190, 2).Value = "09:55:44"
$ws.Cells.Item(190, 3).Value = "09:00"
$ws.Cells.Item(190, 4).Value = "Bathroom"
$ws.Cells.Item(190, 5).Value = "No Motion"
$ws.Cells.Item(190, 6).Value = "Inactive"

$ws.Cells.Item(191, 1).NumberFormat = "@"
$ws.Cells.Item(191, 1).Value = "2026-02-06"
$ws.Cells.Item(191, 2).Value = "09:55:49"
$ws.Cells.Item(191, 3).Value = "09:00"
$ws.Cells.Item(191, 4).Value = "Bathroom"
$ws.Cells.Item(191, 5).Value = "No Motion"
$ws.Cells.Item(191, 6).Value = "Inactive"


# --- Humidity: append 11 row(s) starting at row 99 ---
$ws = $wb.Worksheets.Item("Humidity")

$ws.Cells.Item(99, 1).NumberFormat = "@"
$ws.Cells.Item(99, 5).NumberFormat = "@"
$ws.Cells.Item(99, 1).Value = "2026-02-06"
$ws.Cells.Item(99, 2).Value = "09:54:52"
$ws.Cells.Item(99, 3).Value = "09:00"
$ws.Cells.Item(99, 4).Value = "Bathroom"
$ws.Cells.Item(99, 5).Value = "70.3%"
$ws.Cells.Item(99, 6).Value = "Active"

$ws.Cells.Item(100, 1).NumberFormat = "@"
$ws.Cells.Item(100, 5).NumberFormat = "@"
$ws.Cells.Item(100, 1).Value = "2026-02-06"
$ws.Cells.Item(100, 2).Value = "09:54:57"
$ws.Cells.Item(100, 3).Value = "09:00"
$ws.Cells.Item(100, 4).Value = "Bathroom"
$ws.Cells.Item(100, 5).Value = "70.3%"
$ws.Cells.Item(100, 6).Value = "Active"

$ws.Cells.Item(101, 1).NumberFormat = "@"
$ws.Cells.Item(101, 5).NumberFormat = "@"
$ws.Cells.Item(101, 1).Value = "2026-02-06"
$ws.Cells.Item(101, 2).Value = "09:55:02"
$ws.Cells.Item(101, 3).Value = "09:00"
$ws.Cells.Item(101, 4).Value = "Bathroom"
$ws.Cells.Item(101, 5).Value = "70.4%"
$ws.Cells.Item(101, 6).Value = "Active"

$ws.Cells.Item(102, 1).NumberFormat = "@"
$ws.Cells.Item(102, 5).NumberFormat = "@"
$ws.Cells.Item(102, 1).Value = "2026-02-06"
$ws.Cells.Item(102, 2).Value = "09:55:07"
$ws.Cells.Item(102, 3).Value = "09:00"
$ws.Cells.Item(102, 4).Value = "Bathroom"
$ws.Cells.Item(102, 5).Value = "70.5%"
$ws.Cells.Item(102, 6).Value = "Active"

$ws.Cells.Item(103, 1).NumberFormat = "@"
$ws.Cells.Item(103, 5).NumberFormat = "@"
$ws.Cells.Item(103, 1).Value = "2026-02-06"
$ws.Cells.Item(103, 2).Value = "09:55:12"
$ws.Cells.Item(103, 3).Value = "09:00"
$ws.Cells.Item(103, 4).Value = "Bathroom"
$ws.Cells.Item(103, 5).Value = "70.5%"
$ws.Cells.Item(103, 6).Value = "Active"

$ws.Cells.Item(104, 1).NumberFormat = "@"
$ws.Cells.Item(104, 5).NumberFormat = "@"
$ws.Cells.Item(104, 1).Value = "2026-02-06"
$ws.Cells.Item(104, 2).Value = "09:55:17"
$ws.Cells.Item(104, 3).Value = "09:00"
$ws.Cells.Item(104, 4).Value = "Bathroom"
$ws.Cells.Item(104, 5).Value = "70.6%"
$ws.Cells.Item(104, 6).Value = "Active"

$ws.Cells.Item(105, 1).NumberFormat = "@"
$ws.Cells.Item(105, 5).NumberFormat = "@"
$ws.Cells.Item(105, 1).Value = "2026-02-06"
$ws.Cells.Item(105, 2).Value = "09:55:22"
$ws.Cells.Item(105, 3).Value = "09:00"
$ws.Cells.Item(105, 4).Value = "Bathroom"
$ws.Cells.Item(105, 5).Value = "70.5%"
$ws.Cells.Item(105, 6).Value = "Active"

$ws.Cells.Item(106, 1).NumberFormat = "@"
$ws.Cells.Item(106, 5).NumberFormat = "@"
$ws.Cells.Item(106, 1).Value = "2026-02-06"
$ws.Cells.Item(106, 2).Value = "09:55:27"
$ws.Cells.Item(106, 3).Value = "09:00"
$ws.Cells.Item(106, 4).Value = "Bathroom"
$ws.Cells.Item(106, 5).Value = "70.4%"
$ws.Cells.Item(106, 6).Value = "Active"

$ws.Cells.Item(107, 1).NumberFormat = "@"
$ws.Cells.Item(107, 5).NumberFormat = "@"
$ws.Cells.Item(107, 1).Value = "2026-02-06"
$ws.Cells.Item(107, 2).Value = "09:55:37"
$ws.Cells.Item(107, 3).Value = "09:00"
$ws.Cells.Item(107, 4).Value = "Bathroom"
$ws.Cells.Item(107, 5).Value = "70.5%"
$ws.Cells.Item(107, 6).Value = "Active"

$ws.Cells.Item(108, 1).NumberFormat = "@"
$ws.Cells.Item(108, 5).NumberFormat = "@"
$ws.Cells.Item(108, 1).Value = "2026-02-06"
$ws.Cells.Item(108, 2).Value = "09:55:42"
$ws.Cells.Item(108, 3).Value = "09:00"
$ws.Cells.Item(108, 4).Value = "Bathroom"
$ws.Cells.Item(108, 5).Value = "70.5%"
$ws.Cells.Item(108, 6).Value = "Active"

$ws.Cells.Item(109, 1).NumberFormat = "@"
$ws.Cells.Item(109, 5).NumberFormat = "@"
$ws.Cells.Item(109, 1).Value = "2026-02-06"
$ws.Cells.Item(109, 2).Value = "09:55:47"
$ws.Cells.Item(109, 3).Value = "09:00"
$ws.Cells.Item(109, 4).Value = "Bathroom"
$ws.Cells.Item(109, 5).Value = "70.5%"
$ws.Cells.Item(109, 6).Value = "Active"


# --- Temperature: append 11 row(s) starting at row 99 ---
$ws = $wb.Worksheets.Item("Temperature")

$ws.Cells.Item(99, 1).NumberFormat = "@"
$ws.Cells.Item(99, 1).Value = "2026-02-06"
$ws.Cells.Item(99, 2).Value = "09:54:52"
$ws.Cells.Item(99, 3).Value = "09:00"
$ws.Cells.Item(99, 4).Value = "Bathroom"
$ws.Cells.Item(99, 5).Value = "27.8C"
$ws.Cells.Item(99, 6).Value = "Active"

$ws.Cells.Item(100, 1).NumberFormat = "@"
$ws.Cells.Item(100, 1).Value = "2026-02-06"
$ws.Cells.Item(100, 2).Value = "09:54:57"
$ws.Cells.Item(100, 3).Value = "09:00"
$ws.Cells.Item(100, 4).Value = "Bathroom"
$ws.Cells.Item(100, 5).Value = "27.9C"
$ws.Cells.Item(100, 6).Value = "Active"

$ws.Cells.Item(101, 1).NumberFormat = "@"
$ws.Cells.Item(101, 1).Value = "2026-02-06"
$ws.Cells.Item(101, 2).Value = "09:55:02"
$ws.Cells.Item(101, 3).Value = "09:00"
$ws.Cells.Item(101, 4).Value = "Bathroom"
$ws.Cells.Item(101, 5).Value = "27.8C"
$ws.Cells.Item(101, 6).Value = "Active"

$ws.Cells.Item(102, 1).NumberFormat = "@"
$ws.Cells.Item(102, 1).Value = "2026-02-06"
$ws.Cells.Item(102, 2).Value = "09:55:07"
$ws.Cells.Item(102, 3).Value = "09:00"
$ws.Cells.Item(102, 4).Value = "Bathroom"
$ws.Cells.Item(102, 5).Value = "27.8C"
$ws.Cells.Item(102, 6).Value = "Active"

$ws.Cells.Item(103, 1).NumberFormat = "@"
$ws.Cells.Item(103, 1).Value = "2026-02-06"
$ws.Cells.Item(103, 2).Value = "09:55:12"
$ws.Cells.Item(103, 3).Value = "09:00"
$ws.Cells.Item(103, 4).Value = "Bathroom"
$ws.Cells.Item(103, 5).Value = "27.8C"
$ws.Cells.Item(103, 6).Value = "Active"

$ws.Cells.Item(104, 1).NumberFormat = "@"
$ws.Cells.Item(104, 1).Value = "2026-02-06"
$ws.Cells.Item(104, 2).Value = "09:55:17"
$ws.Cells.Item(104, 3).Value = "09:00"
$ws.Cells.Item(104, 4).Value = "Bathroom"
$ws.Cells.Item(104, 5).Value = "27.9C"
$ws.Cells.Item(104, 6).Value = "Active"

$ws.Cells.Item(105, 1).NumberFormat = "@"
$ws.Cells.Item(105, 1).Value = "2026-02-06"
$ws.Cells.Item(105, 2).Value = "09:55:23"
$ws.Cells.Item(105, 3).Value = "09:00"
$ws.Cells.Item(105, 4).Value = "Bathroom"
$ws.Cells.Item(105, 5).Value = "27.9C"
$ws.Cells.Item(105, 6).Value = "Active"

$ws.Cells.Item(106, 1).NumberFormat = "@"
$ws.Cells.Item(106, 1).Value = "2026-02-06"
$ws.Cells.Item(106, 2).Value = "09:55:28"
$ws.Cells.Item(106, 3).Value = "09:00"
$ws.Cells.Item(106, 4).Value = "Bathroom"
$ws.Cells.Item(106, 5).Value = "27.8C"
$ws.Cells.Item(106, 6).Value = "Active"

$ws.Cells.Item(107, 1).NumberFormat = "@"
$ws.Cells.Item(107, 1).Value = "2026-02-06"
$ws.Cells.Item(107, 2).Value = "09:55:38"
$ws.Cells.Item(107, 3).Value = "09:00"
$ws.Cells.Item(107, 4).Value = "Bathroom"
$ws.Cells.Item(107, 5).Value = "27.8C"
$ws.Cells.Item(107, 6).Value = "Active"

$ws.Cells.Item(108, 1).NumberFormat = "@"
$ws.Cells.Item(108, 1).Value = "2026-02-06"
$ws.Cells.Item(108, 2).Value = "09:55:43"
$ws.Cells.Item(108, 3).Value = "09:00"
$ws.Cells.Item(108, 4).Value = "Bathroom"
$ws.Cells.Item(108, 5).Value = "27.9C"
$ws.Cells.Item(108, 6).Value = "Active"

$ws.Cells.Item(109, 1).NumberFormat = "@"
$ws.Cells.Item(109, 1).Value = "2026-02-06"
$ws.Cells.Item(109, 2).Value = "09:55:48"
$ws.Cells.Item(109, 3).Value = "09:00"
$ws.Cells.Item(109, 4).Value = "Bathroom"
$ws.Cells.Item(109, 5).Value = "27.8C"
$ws.Cells.Item(109, 6).Value = "Active"


# --- Proximity: append 2 row(s) starting at row 9 ---
$ws = $wb.Worksheets.Item("Proximity")

$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "2026-02-06"
$ws.Cells.Item(9, 2).Value = "09:55:15"
$ws.Cells.Item(9, 3).Value = "09:00"
$ws.Cells.Item(9, 4).Value = "Bathroom Door"
$ws.Cells.Item(9, 5).Value = "EXIT"
$ws.Cells.Item(9, 6).Value = "User EXITED Bathroom"

$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "2026-02-06"
$ws.Cells.Item(10, 2).Value = "09:55:23"
$ws.Cells.Item(10, 3).Value = "09:00"
$ws.Cells.Item(10, 4).Value = "Bathroom Door"
$ws.Cells.Item(10, 5).Value = "ENTER"
$ws.Cells.Item(10, 6).Value = "User ENTERED Bathroom"

